$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("A1").Value = "Plate_Number"
$ws.Range("B1").Value = "Owner_Name"
$ws.Range("C1").Value = "Vehicle_Type"
$ws.Range("D1").Value = "Registration_Date"

# Full data set for rows 2-11 (plate, owner, vehicle type, registration date)
$data = @(
    @("ABC 123 XYZ", "John Doe", "Sedan", "2024-01-15"),
    @("DEF 456 UVW", "Jane Smith", "SUV", "2024-02-20"),
    @("GHI 789 RST", "Bob Johnson", "Truck", "2024-03-10"),
    @("JKL 012 PQR", "Alice Williams", "Sedan", "2024-03-25"),
    @("MNO 345 LMN", "Charlie Brown", "Hatchback", "2024-04-05"),
    @("PQR 678 JKL", "Diana Davis", "SUV", "2024-04-18"),
    @("STU 901 GHI", "Eva Wilson", "Sedan", "2024-05-02"),
    @("VWX 234 DEF", "Frank Miller", "Truck", "2024-05-15"),
    @("YZA 567 ABC", "Grace Taylor", "SUV", "2024-06-01"),
    @("BCD 890 ZYX", "Henry Anderson", "Sedan", "2024-06-20")
)

# Format the Registration_Date column as text first so Excel does not
# auto-convert the "YYYY-MM-DD" strings into date serial numbers.
$ws.Range("D2:D11").NumberFormat = "@"

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $row++
}

# Restore default (unstyled) formatting for the data rows so the cells
# keep the plain inline-string representation rather than a custom
# number format style.
$ws.Range("D2:D11").Style = "Normal"
